# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement table (rows 16-74, columns E "Periodo Mora" and
# F "Valor Mora") is re-sorted from newest-period-first to oldest-period-
# first (chronological ascending), while each period keeps the "Valor
# Mora" value it already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the 59 period codes in chronological order: 1607 .. 1612, 1701 .. 1712,
# 1801 .. 1812, 1901 .. 1912, 2001 .. 2012, 2101 .. 2105
$periods = New-Object System.Collections.ArrayList
for ($yy = 16; $yy -le 21; $yy++) {
    $startMonth = 1
    $endMonth = 12
    if ($yy -eq 16) { $startMonth = 7 }
    if ($yy -eq 21) { $endMonth = 5 }
    for ($mm = $startMonth; $mm -le $endMonth; $mm++) {
        $code = "{0:D2}{1:D2}" -f $yy, $mm
        [void]$periods.Add($code)
    }
}

# The "Valor Mora" associated with each period (unchanged per-period,
# only the row each period lands on changes): the 26 oldest periods
# (1607-1808) carry 27578, the next 32 (1809-2104) carry 31249, and the
# newest period (2105) carries 24999.
$values = New-Object System.Collections.ArrayList
for ($i = 0; $i -lt 26; $i++) { [void]$values.Add(27578) }
for ($i = 0; $i -lt 32; $i++) { [void]$values.Add(31249) }
[void]$values.Add(24999)

$firstRow = 16
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
